# Add data for 2022-08-31
# - Rename sheet "Through 2022-08-22" -> "Through 2022-08-23"
# - Update the header label in I1 from "2022 (through 08-22)" to "2022 (through 08-23)"
# - Update September value I9: 129 -> 130
# - Update Total value I14: 1100 -> 1101

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet (tab name / workbook.xml sheet name)
$ws.Name = "Through 2022-08-23"

# Update the shared string used as the column header in I1
$ws.Range("I1").Value = "2022 (through 08-23)"

# Update the September row value (row 9, column I)
$ws.Range("I9").Value = 130

# Update the Total row value (row 14, column I)
$ws.Range("I14").Value = 1101
